$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the diff (cryptos price/volume refresh,
# including a rank swap between Maker and VeChain at rows 43-44).
$ws.Range("D2").Value = "61.339.50"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.929.62"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.40"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.18"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "3.414.49"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "61.317.45"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "2.928.10"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.82"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.04"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.48"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.81"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.61"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.48"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.06"
$ws.Range("E41").Value = "  +5.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.278"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.698.87"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0343"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.47"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "362.38"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.50"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -0.11%  "
